# "Minor update: Just checking whether this works with .xlsx"
#
# Actual content changes reconstructed from the OOXML diff:
#   1. The "Knoppen" sheet gets a new column inserted after "Actie" (col A):
#        - new column B = "activeringsEvent" (event that triggers the action)
#        - the old "conditie" column is relabeled "activeringsConditie" and
#          shifts from B to C
#        - beschrijving/toelichting/Rollen shift from C/D/E to D/E/F
#      A couple of cells in the new layout get real new content:
#        B2 = "Event"      (the "type" row gains a type for the new column)
#        E2 = "TEKST"      (the type row also now reaches column E)
#        B3 = "Klik op knop" (the activeringsEvent for "Vernieuw vanuit Spin")
#   2. "Knoppen" becomes the active sheet/tab (it was "Tabbladen" before).
#   3. Cosmetic-only column width / row height tweaks on the "Knoppen" sheet
#      are approximated as closely as the host lets us.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Knoppen")

# --- insert the new "activeringsEvent" column before the old "conditie" column (B) ---
$ws.Columns.Item(2).Insert()

# --- header row ---
$ws.Range("B1").Value = "activeringsEvent"
$ws.Range("C1").Value = "activeringsConditie"

# --- "type" row ---
$ws.Range("B2").Value = "Event"
$ws.Range("E2").Value = "TEKST"

# --- first data row: new activeringsEvent value ---
$ws.Range("B3").Value = "Klik op knop"

# --- approximate column widths to match the widened / re-fitted layout ---
$ws.Columns.Item(1).ColumnWidth = 13.8776041666667
$ws.Columns.Item(2).ColumnWidth = 19.5924479166667
$ws.Columns.Item(3).ColumnWidth = 31.5924479166667
$ws.Columns.Item(4).ColumnWidth = 36.5924479166667
$ws.Columns.Item(5).ColumnWidth = 57.7369791666667

# --- approximate row heights (re-fit after the column/content changes) ---
$ws.Rows.Item(3).RowHeight = 90
$ws.Rows.Item(4).RowHeight = 225
$ws.Rows.Item(5).RowHeight = 60
$ws.Rows.Item(6).RowHeight = 135

# --- "Knoppen" becomes the active/selected sheet ---
$ws.Activate()
$ws.Range("D3").Select()
